# Fixed filter issue for products with multiple type tags
# -> Add the missing "Aveda" product row (it uses the paginated pg1/pg2/pg3
#    "Line" tagging scheme like As I Am / Jessicurl / Kinky Curly / Ouidad /
#    Dippity Do, rather than a single Line tag), so the Amazon-search-based
#    filter correctly picks up all of its pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row right after the current last row (row 34)
$newRow = 35
$ws.Cells.Item($newRow, 1).Value = "Aveda"
$ws.Cells.Item($newRow, 2).Value = "pg1"
$ws.Cells.Item($newRow, 3).Value = "https://www.amazon.com/s?k=Aveda&rh=n%3A3760911%2Cp_89%3AAveda&dc&qid=1650933740&rnid=2528832011&ref=sr_nr_p_89_1"

# Scroll the view up a bit and select the newly-added cell, matching the
# author's on-screen state after the edit (topLeftCell A8, selection C35)
$excel.Goto($ws.Range("A8"), $true)
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select()
